$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a weekly log of "Coliflor" (cauliflower) price records for
# "Macroferia Regional de Talca". The edit inserts two brand-new weekly
# records into the middle of the table:
#   - one new record at row 126 (pushing the former rows 126-185 down to 127-186)
#   - one new record at row 186, once the sheet has already shifted once
#     (the former row 185 now sits at row 186, so inserting there pushes the
#     former rows 186-191 -- now at 187-192 -- down to 188-193)
# All rows below each insertion point shift down by one, and the sheet
# dimension grows from A1:R191 to A1:R193.
# ---------------------------------------------------------------------------

function Add-ColiflorRow {
    param($RowIndex, $FechaSerial, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Origen, $PrecioKg)

    # Insert a new blank row, shifting the current row (and everything below) down
    $ws.Rows.Item($RowIndex).Insert()

    $ws.Cells.Item($RowIndex, 1).Value = 5
    $ws.Cells.Item($RowIndex, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($RowIndex, 3).Value = "Maule"

    $dCell = $ws.Cells.Item($RowIndex, 4)
    $dCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $dCell.Value = $FechaSerial

    $ws.Cells.Item($RowIndex, 5).Value = 7
    $ws.Cells.Item($RowIndex, 6).Value = 100112008
    $ws.Cells.Item($RowIndex, 7).Value = "Coliflor"
    $ws.Cells.Item($RowIndex, 8).Value = "Sin especificar"
    $ws.Cells.Item($RowIndex, 9).Value = $Calidad
    $ws.Cells.Item($RowIndex, 10).Value = $Volumen
    $ws.Cells.Item($RowIndex, 11).Value = $PrecioMin
    $ws.Cells.Item($RowIndex, 12).Value = $PrecioMax
    $ws.Cells.Item($RowIndex, 13).Value = $PrecioProm
    $ws.Cells.Item($RowIndex, 14).Value = "`$/unidad"
    $ws.Cells.Item($RowIndex, 15).Value = $Origen
    $ws.Cells.Item($RowIndex, 16).Value = $PrecioKg
    $ws.Cells.Item($RowIndex, 17).Value = 1
    $ws.Cells.Item($RowIndex, 18).Value = "Hortaliza"
}

# New record inserted at row 126 (date serial 44567 = 2022-01-06)
Add-ColiflorRow 126 44567 "Primera" 2000 800 800 800 "Región del Maule" 800

# New record inserted at row 186 (date serial 44568 = 2022-01-07); by this point
# the previous insert has already shifted the old row 185 down to row 186, so
# inserting here pushes that row (and everything below it) down by one more.
Add-ColiflorRow 186 44568 "Primera" 2000 800 800 800 "Región del Maule" 800
